$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.696.67"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "2.469.66"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("D5").Value = "'323.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "'104.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'36.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "'18.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.61%  "
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "2.853.83"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "2.449.16"
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "46.607.36"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").Value = "'12.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "'6.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").Value = "'70.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'2.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.05%  "
$ws.Range("D24").Value = "'249.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").Value = "'26.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'9.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("D31").Value = "'0.133"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.82%  "
$ws.Range("D32").Value = "'49.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").Value = "'19.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "'0.0767"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "'2.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'123.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").Value = "'20.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").Value = "1.977.67"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "'2.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").Value = "'1.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("E49").Value = "  +16.14%  "
$ws.Range("D50").Value = "'8.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").Value = "'79.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.24%  "
